$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to remain text (product codes look numeric but are stored as text)
$ws.Range("A2:A28").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '1008277'
$ws.Range("B2").Value = 'Personalised 80th Birthday Book ''Memory Lane'''
$ws.Range("D2").Value = 'lucysworld'
$ws.Range("F2").Value = 'https://www.notonthehighstreet.com/lucysworld/product/personalised-80th-birthday-book-memory-lane'
$ws.Range("G2").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1008277&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 3
$ws.Range("A3").Value = '1062346'
$ws.Range("B3").Value = 'Tree Of Life Rose Gold Foil Scarf Gift'
$ws.Range("D3").Value = 'studiohop'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 'https://www.notonthehighstreet.com/studiohop/product/tree-of-life-foil-scarf-letterbox-gift'
$ws.Range("G3").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1062346&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 4
$ws.Range("A4").Value = '1130272'
$ws.Range("B4").Value = 'Personalised Family Seaside Beach Pebble Picture'
$ws.Range("D4").Value = 'ladedaliving'
$ws.Range("F4").Value = 'https://www.notonthehighstreet.com/ladedaliving/product/personalised-family-seaside-beach-pebble-picture'
$ws.Range("G4").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1130272&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 5
$ws.Range("A5").Value = '1179458'
$ws.Range("B5").Value = 'Long Wildflower Amazing Friend Trinket Dish'
$ws.Range("D5").Value = 'lisaangeljewellery'
$ws.Range("F5").Value = 'https://www.notonthehighstreet.com/lisaangeljewellery/product/long-wildflower-amazing-friend-trinket-dish'
$ws.Range("G5").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1179458&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 6
$ws.Range("A6").Value = '1218093'
$ws.Range("B6").Value = 'Birth Month Flower Earring Studs'
$ws.Range("D6").Value = 'attic'
$ws.Range("F6").Value = 'https://www.notonthehighstreet.com/attic/product/birth-month-flower-earring-studs'
$ws.Range("G6").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1218093&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 8
$ws.Range("A8").Value = '1323961'
$ws.Range("B8").Value = 'Muslin Swaddle Blanket Hello World Newborn Baby Shower Gift'
$ws.Range("D8").Value = 'geople'
$ws.Range("F8").Value = 'https://www.notonthehighstreet.com/geople/product/muslin-swaddle-blanket-sunshine-newborn-baby-gift'
$ws.Range("G8").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1323961&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 9
$ws.Range("A9").Value = '133151'
$ws.Range("B9").Value = 'Sheep Sofa Tidy, Remote Control Holder'
$ws.Range("D9").Value = 'jomanda'
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 'https://www.notonthehighstreet.com/jomanda/product/sheepey-sofa-tidy-coffe'
$ws.Range("G9").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=133151&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 10
$ws.Range("A10").Value = '1348765'
$ws.Range("B10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1348765&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 11
$ws.Range("A11").Value = '1353406'
$ws.Range("G11").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1353406&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 12
$ws.Range("A12").Value = '1359483'
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1359483&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 13
$ws.Range("A13").Value = '1385338'
$ws.Range("G13").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1385338&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 14
$ws.Range("A14").Value = '1400456'
$ws.Range("G14").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1400456&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 15
$ws.Range("A15").Value = '1462617'
$ws.Range("G15").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1462617&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 16
$ws.Range("A16").Value = '1476087'
$ws.Range("G16").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=1476087&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 17
$ws.Range("A17").Value = '350209'
$ws.Range("G17").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=350209&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 18
$ws.Range("A18").Value = '412675'
$ws.Range("G18").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=412675&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 19
$ws.Range("A19").Value = '469358'
$ws.Range("G19").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=469358&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 20
$ws.Range("A20").Value = '610619'
$ws.Range("E20").Value = 2
$ws.Range("G20").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=610619&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 21
$ws.Range("A21").Value = '631040'
$ws.Range("E21").Value = 2
$ws.Range("G21").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=631040&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 22
$ws.Range("A22").Value = '748820'
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=748820&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 23
$ws.Range("A23").Value = '764151'
$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=764151&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 24
$ws.Range("A24").Value = '786481'
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=786481&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 25
$ws.Range("A25").Value = '864619'
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=864619&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 26
$ws.Range("A26").Value = '876141'
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=876141&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 27
$ws.Range("A27").Value = '905169'
$ws.Range("E27").Value = 2
$ws.Range("G27").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=905169&displayFeedbackType=PRODUCT&timeFrame=ALL'
# Row 28
$ws.Range("A28").Value = '999592'
$ws.Range("E28").Value = 2
$ws.Range("G28").Value = 'https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=999592&displayFeedbackType=PRODUCT&timeFrame=ALL'
